$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Factor" / "Def Vals" header block for the white-LED section (rows 8-11)
$ws.Range("K8").Value = "Factpr"
$ws.Range("L8").Value = 0.5
$ws.Range("N8").Value = "Def Vals"

# Preserve the original 0-255 default values in the new N column
$ws.Range("N9").Value = 104
$ws.Range("N10").Value = 135
$ws.Range("N11").Value = 130

# Drive the B column (0-255 input) off the new default column * factor
$ws.Range("B9").Formula = '=N9*$L$8'
$ws.Range("B10").Formula = '=N10*$L$8'
$ws.Range("B11").Formula = '=N11*$L$8'

# Restore the original selection used by the authored workbook
$ws.Range("L9").Select()
